# The sheet contains a long historical price table (Betarraga / Femacal de
# La Calera). This edit inserts one new reporting period (2 rows: "Primera"
# and "Segunda" quality grades) right before the existing row 1033, pushing
# all the old rows 1033:1161 down to 1035:1163 (dimension grows from
# A1:R1161 to A1:R1163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1033/1034 - everything below (old 1033:1161)
# shifts down to 1035:1163, carrying its data/styles with it.
$ws.Rows("1033:1034").Insert()

# --- New row 1033 (Calidad = Primera) ---
$ws.Cells.Item(1033, 1).Value = 3
$ws.Cells.Item(1033, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1033, 3).Value = "Coquimbo"
$ws.Cells.Item(1033, 4).Value = 45124
$ws.Cells.Item(1033, 5).Value = 5
$ws.Cells.Item(1033, 6).Value = 100114014
$ws.Cells.Item(1033, 7).Value = "Betarraga"
$ws.Cells.Item(1033, 8).Value = "Sin especificar"
$ws.Cells.Item(1033, 9).Value = "Primera"
$ws.Cells.Item(1033, 10).Value = 3500
$ws.Cells.Item(1033, 11).Value = 550
$ws.Cells.Item(1033, 12).Value = 600
$ws.Cells.Item(1033, 13).Value = 573
$ws.Cells.Item(1033, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(1033, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1033, 16).Value = 143
$ws.Cells.Item(1033, 17).Value = 4
$ws.Cells.Item(1033, 18).Value = "Hortaliza"

# --- New row 1034 (Calidad = Segunda) ---
$ws.Cells.Item(1034, 1).Value = 3
$ws.Cells.Item(1034, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1034, 3).Value = "Coquimbo"
$ws.Cells.Item(1034, 4).Value = 45124
$ws.Cells.Item(1034, 5).Value = 5
$ws.Cells.Item(1034, 6).Value = 100114014
$ws.Cells.Item(1034, 7).Value = "Betarraga"
$ws.Cells.Item(1034, 8).Value = "Sin especificar"
$ws.Cells.Item(1034, 9).Value = "Segunda"
$ws.Cells.Item(1034, 10).Value = 1900
$ws.Cells.Item(1034, 11).Value = 450
$ws.Cells.Item(1034, 12).Value = 450
$ws.Cells.Item(1034, 13).Value = 450
$ws.Cells.Item(1034, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(1034, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1034, 16).Value = 112
$ws.Cells.Item(1034, 17).Value = 4
$ws.Cells.Item(1034, 18).Value = "Hortaliza"
